$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRange, $newValue)
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $newValue
    $cellRange.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "34.671.48"
Set-TextValue $ws.Range("E2") "  +0.55%  "
Set-TextValue $ws.Range("D3") "1.828.83"
Set-TextValue $ws.Range("E3") "  +1.08%  "
Set-TextValue $ws.Range("D5") "227.24"
Set-TextValue $ws.Range("E5") "  +0.64%  "
Set-TextValue $ws.Range("D6") "0.611"
Set-TextValue $ws.Range("E6") "  +1.69%  "
Set-TextValue $ws.Range("E7") "  +0.22%  "
Set-TextValue $ws.Range("D8") "44.03"
Set-TextValue $ws.Range("E8") "  +21.16%  "
Set-TextValue $ws.Range("D9") "0.300"
Set-TextValue $ws.Range("E9") "  +2.06%  "
Set-TextValue $ws.Range("D10") "0.0684"
Set-TextValue $ws.Range("E10") "  +0.31%  "
Set-TextValue $ws.Range("E11") "  +3.91%  "
Set-TextValue $ws.Range("D12") "2.092.56"
Set-TextValue $ws.Range("E12") "  +1.08%  "
Set-TextValue $ws.Range("B13") "Chainlink"
Set-TextValue $ws.Range("C13") "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D13") "11.23"
Set-TextValue $ws.Range("E13") "  -1.26%  "
Set-TextValue $ws.Range("B14") "WrappedEther"
Set-TextValue $ws.Range("C14") "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D14") "1.823.11"
Set-TextValue $ws.Range("E14") "  +0.80%  "
Set-TextValue $ws.Range("D15") "4.70"
Set-TextValue $ws.Range("E15") "  +5.77%  "
Set-TextValue $ws.Range("D16") "0.650"
Set-TextValue $ws.Range("E16") "  +2.98%  "
Set-TextValue $ws.Range("D17") "34.645.06"
Set-TextValue $ws.Range("E17") "  +0.60%  "
Set-TextValue $ws.Range("D18") "68.41"
Set-TextValue $ws.Range("E18") "  -0.34%  "
Set-TextValue $ws.Range("D19") "242.73"
Set-TextValue $ws.Range("E19") "  -0.17%  "
Set-TextValue $ws.Range("D20") "0.0₃0786"
Set-TextValue $ws.Range("E20") "  +1.46%  "
Set-TextValue $ws.Range("D21") "11.97"
Set-TextValue $ws.Range("E21") "  +6.11%  "
Set-TextValue $ws.Range("D22") "4.68"
Set-TextValue $ws.Range("E22") "  +13.97%  "
Set-TextValue $ws.Range("E23") "  +0.24%  "
Set-TextValue $ws.Range("D24") "2.19"
Set-TextValue $ws.Range("E24") "  -1.34%  "
Set-TextValue $ws.Range("D25") "171.39"
Set-TextValue $ws.Range("E25") "  +0.03%  "
Set-TextValue $ws.Range("D26") "7.93"
Set-TextValue $ws.Range("E26") "  +0.40%  "
Set-TextValue $ws.Range("D27") "17.86"
Set-TextValue $ws.Range("E27") "  +3.00%  "
Set-TextValue $ws.Range("D28") "0.122"
Set-TextValue $ws.Range("E28") "  +0.82%  "
Set-TextValue $ws.Range("E29") "  +0.33%  "
Set-TextValue $ws.Range("D30") "3.89"
Set-TextValue $ws.Range("E30") "  +1.95%  "
Set-TextValue $ws.Range("D31") "1.25"
Set-TextValue $ws.Range("E31") "  +1.85%  "
Set-TextValue $ws.Range("D32") "3.97"
Set-TextValue $ws.Range("E32") "  +1.12%  "
Set-TextValue $ws.Range("D33") "0.0524"
Set-TextValue $ws.Range("E33") "  +1.30%  "
Set-TextValue $ws.Range("D34") "1.85"
Set-TextValue $ws.Range("E34") "  +2.87%  "
Set-TextValue $ws.Range("D35") "90.81"
Set-TextValue $ws.Range("E35") "  +11.82%  "
Set-TextValue $ws.Range("D36") "0.663"
Set-TextValue $ws.Range("E36") "  +1.31%  "
Set-TextValue $ws.Range("D37") "15.43"
Set-TextValue $ws.Range("E37") "  +14.86%  "
Set-TextValue $ws.Range("D38") "1.329.67"
Set-TextValue $ws.Range("E38") "  -2.39%  "
Set-TextValue $ws.Range("B39") "RenderToken"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D39") "2.43"
Set-TextValue $ws.Range("E39") "  +2.10%  "
Set-TextValue $ws.Range("B40") "TrustWalletToken"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D40") "1.06"
Set-TextValue $ws.Range("E40") "  +0.30%  "
Set-TextValue $ws.Range("D41") "0.0192"
Set-TextValue $ws.Range("E41") "  +2.71%  "
Set-TextValue $ws.Range("D42") "0.968"
Set-TextValue $ws.Range("E42") "  +3.06%  "
Set-TextValue $ws.Range("D43") "1.23"
Set-TextValue $ws.Range("E43") "  +5.62%  "
Set-TextValue $ws.Range("B44") "MXToken"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D44") "2.82"
Set-TextValue $ws.Range("E44") "  +1.07%  "
Set-TextValue $ws.Range("B45") "HuobiToken"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D45") "2.43"
Set-TextValue $ws.Range("E45") "  +0.27%  "
Set-TextValue $ws.Range("E46") "  +3.72%  "
Set-TextValue $ws.Range("D47") "1.991.18"
Set-TextValue $ws.Range("E47") "  +1.06%  "
Set-TextValue $ws.Range("D48") "5.94"
Set-TextValue $ws.Range("E48") "  +1.79%  "
Set-TextValue $ws.Range("E49") "  +0.32%  "
Set-TextValue $ws.Range("D50") "101.88"
Set-TextValue $ws.Range("E50") "  -0.76%  "
Set-TextValue $ws.Range("E51") "  +18.97%  "
